$wb = $excel.ActiveWorkbook

$wsCoef = $wb.Worksheets.Item("Coefficients")
$wsCoef.Range("A5").Value = "MonthOctober"
$wsCoef.Range("A7").Value = "DRM:MonthOctober"
$wsCoef.Range("A8").Value = "Depth:MonthOctober"

$wsStats = $wb.Worksheets.Item("Fullmodel_statistics")
$wsStats.Range("A2").Value = 0.7929473038810525
$wsStats.Range("B2").Value = 0.7451659124689877
$wsStats.Range("E2").Value = 0.00000008808597263347452
